# Advocacy / drawdown calcs.xlsx
# "added increased pumping(x1.25) for initial baseline"
#
# Duplicate the existing Ensemble/Number/Initial-Drawdown-Aq baseline table
# (rows 2-12) down to rows 18-28, then add a new "more pump drawdown" column
# (D) with the increased-pumping drawdown values and a "% Diff" column (E)
# comparing the two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the original header+data block (A2:C12) down to A18:C28.
#    This reproduces the per-cell styles (borders/fonts/number formats)
#    that the source block already uses, exactly as a copy/paste in Excel
#    would.
$src = $ws.Range("A2:C12")
$dst = $ws.Range("A18")
$src.Copy($dst)

# 2) New header cell for the extra column.
$ws.Range("D18").Value = "more pump drawdown"

# 3) First three ensembles have no comparable "more pump" figure yet.
$ws.Range("D19").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("D21").Value = "-"

# 4) Increased-pumping (~x1.25) drawdown readings for the remaining
#    ensembles, plus the %-difference formula against column C.
$ws.Range("D22").Value = -10.02
$ws.Range("D23").Value = -2.19
$ws.Range("D24").Value = -2.21
$ws.Range("D25").Value = -2.21
$ws.Range("D26").Value = -1.015
$ws.Range("D27").Value = -1.44
$ws.Range("D28").Value = -10.21

$ws.Range("E22").Formula = "=(1-(C22/D22))*100"
$ws.Range("E23").Formula = "=(1-(C23/D23))*100"
$ws.Range("E24").Formula = "=(1-(C24/D24))*100"
$ws.Range("E25").Formula = "=(1-(C25/D25))*100"
$ws.Range("E26").Formula = "=(1-(C26/D26))*100"
$ws.Range("E27").Formula = "=(1-(C27/D27))*100"
$ws.Range("E28").Formula = "=(1-(C28/D28))*100"

# 5) Leave the selection where the author ended up.
$null = $ws.Range("D28").Select()
